$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "RaFwC907"
$ws.Range("B2").Value = 231031243
$ws.Range("C2").Value = "fntqnvt71"
$ws.Range("D2").Value = "R8Em#2&j"
$ws.Range("F2").Value = "DtZDnUeF"
$ws.Range("G2").Value = "defP"
